# Update column G ("K") values for rows 3-25 on the active worksheet.
# These values were regenerated (K computed from source data instead of the
# previous "Strike#" based value) per the commit message:
#   "regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    3  = 2
    4  = 2
    5  = 0
    6  = 0
    7  = 0
    8  = 0
    9  = 1
    10 = 1
    11 = 1
    12 = 0
    13 = 2
    14 = 2
    15 = 1
    16 = 2
    17 = 0
    18 = 2
    19 = 0
    20 = 0
    21 = 2
    22 = 1
    23 = 1
    24 = 2
    25 = 0
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
